$d = $word.ActiveDocument
$d.Content.Find.Execute("REQUERIMENTO", $true, $false, $false, $false, $false,
                         $true, 1, $false, "REQUERIMENTO", 2)
